$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing value for the existing "Day 4" row
$ws.Range("B7").Value = 7

# Add new rows for Day 5 .. Day 10
$ws.Range("A8").Value = "Day 5"
$ws.Range("B8").Value = 1

$ws.Range("A9").Value = "Day 6"
$ws.Range("B9").Value = 1

$ws.Range("A10").Value = "Day 7"
$ws.Range("B10").Value = 1

$ws.Range("A11").Value = "Day 8 "
$ws.Range("B11").Value = 1

$ws.Range("A12").Value = "Day 9"
$ws.Range("B12").Value = 1

$ws.Range("A13").Value = "Day 10 "
$ws.Range("B13").Value = 2

$ws.Range("C13").Select()
